$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("B2").Value = "NA"
$ws.Range("C2").Value = "NA"
$ws.Range("D2").Value = "NA"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1

# Update the selected cell to G2 (as reflected in the saved view state)
$ws.Range("G2").Select()
